# Update BOM descriptions (Description column, B) and footprint names
# (Footprint column, D) to match the re-uploaded workbook revision, while
# preserving the original cell formatting (style index) of each cell.
#
# Setting .Value directly would make the host normalize away the cell's
# quote-prefix formatting flag (switching its style index). To avoid that,
# after writing the new value we copy the *format only* from an untouched
# neighboring cell in the same column (which still carries the original
# style) back onto the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellTextKeepFormat {
    param($addr, $text, $refAddr)
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $ref = $ws.Range($refAddr)
    $ref.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# Capacitors: distinguish non-polarized (ceramic) vs polarized (electrolytic/tantalum)
Set-CellTextKeepFormat "B2" "贴片无极性电容" "B16"
Set-CellTextKeepFormat "B3" "贴片无极性电容" "B16"
Set-CellTextKeepFormat "B4" "贴片无极性电容" "B16"
Set-CellTextKeepFormat "B5" "贴片无极性电容" "B16"
Set-CellTextKeepFormat "B6" "贴片有极性电容" "B16"
Set-CellTextKeepFormat "B7" "贴片有极性电容" "B16"
Set-CellTextKeepFormat "B8" "贴片有极性电容" "B16"

# LED footprints simplified (drop color/suffix qualifiers)
Set-CellTextKeepFormat "D17" "LED-0603" "D16"
Set-CellTextKeepFormat "D18" "LED-0603" "D16"

# MCU description + footprint suffix cleanup
Set-CellTextKeepFormat "B28" "增强型51单片机" "B16"
Set-CellTextKeepFormat "D28" "LQFP48 7X7" "D16"

# USB-to-serial chip description clarified
Set-CellTextKeepFormat "B29" "USB转串口芯片，免外部晶振，超小封装" "B16"

# Sensor descriptions translated to Chinese
Set-CellTextKeepFormat "B32" "三轴磁力计" "B16"
Set-CellTextKeepFormat "B33" "实时时钟" "B16"
Set-CellTextKeepFormat "B35" "温度湿度压强一体式传感器" "B16"

# Restore the last-saved selection shown in the sheet view
$ws.Range("D40").Select()
